# Auto-generated edit script.
# Source change: reorder monthly rows (Oct/Nov/Dec moved to the front of each
# calendar-year block) and drop column F ("...出口交货值") from Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final data for A2:E63, in the new row order, built from the source values.
$data = New-Object 'object[,]' 62,5

$data[0,0] = "2012-10"; $data[0,1] = [double]-64.39; $data[0,2] = [double]0.01; $data[0,3] = [double]0.13; $data[0,4] = [double]-87.03
$data[1,0] = "2012-11"; $data[1,1] = [double]-72.93000000000001; $data[1,2] = [double]0.01; $data[1,3] = [double]0.14; $data[1,4] = [double]-86.3
$data[2,0] = "2012-12"; $data[2,1] = [double]-80.7; $data[2,2] = [double]0.01; $data[2,3] = [double]0.16; $data[2,4] = [double]-85.91
$data[3,0] = "2012-02"; $data[3,1] = [double]-69.78; $data[3,2] = [double]0.06; $data[3,3] = [double]0.23; $data[3,4] = [double]-34.58
$data[4,0] = "2012-03"; $data[4,1] = [double]-100; $data[4,2] = [double]0; $data[4,3] = [double]0; $data[4,4] = [double]-100
$data[5,0] = "2012-04"; $data[5,1] = [double]160.22; $data[5,2] = [double]0.17; $data[5,3] = [double]0.43; $data[5,4] = [double]-26.28
$data[6,0] = "2012-05"; $data[6,1] = [double]-65.92; $data[6,2] = [double]0.02; $data[6,3] = [double]0.05; $data[6,4] = [double]-91.76000000000001
$data[7,0] = "2012-06"; $data[7,1] = [double]-80.83; $data[7,2] = [double]0.02; $data[7,3] = [double]0.07000000000000001; $data[7,4] = [double]-89.81
$data[8,0] = "2012-07"; $data[8,1] = [double]-92.97; $data[8,2] = [double]0.01; $data[8,3] = [double]0.08; $data[8,4] = [double]-90.45
$data[9,0] = "2012-08"; $data[9,1] = [double]47.25; $data[9,2] = [double]0.11; $data[9,3] = [double]0.19; $data[9,4] = [double]-79.54000000000001
$data[10,0] = "2012-09"; $data[10,1] = [double]-74.89; $data[10,2] = [double]0.02; $data[10,3] = [double]0.38; $data[10,4] = [double]-63.24
$data[11,0] = "2013-10"; $data[11,1] = [double]247.7; $data[11,2] = [double]0.1; $data[11,3] = [double]0.5; $data[11,4] = [double]294.5
$data[12,0] = "2013-11"; $data[12,1] = [double]568; $data[12,2] = [double]0.1; $data[12,3] = [double]0.6; $data[12,4] = [double]322.7
$data[13,0] = "2013-12"; $data[13,1] = [double]-86.40000000000001; $data[13,2] = ""; $data[13,3] = ""; $data[13,4] = [double]-49.3
$data[14,0] = "2013-02"; $data[14,1] = [double]321.5; $data[14,2] = [double]0.3; $data[14,3] = [double]0.5; $data[14,4] = [double]131.6
$data[15,0] = "2013-03"; $data[15,1] = [double]-100; $data[15,2] = [double]0.2; $data[15,3] = [double]0.2; $data[15,4] = [double]-100
$data[16,0] = "2013-04"; $data[16,1] = [double]55.6; $data[16,2] = [double]0.3; $data[16,3] = [double]0.6; $data[16,4] = [double]35.4
$data[17,0] = "2013-05"; $data[17,1] = [double]13.1; $data[17,2] = [double]0; $data[17,3] = [double]0.2; $data[17,4] = [double]347.5
$data[18,0] = "2013-06"; $data[18,1] = [double]214.5; $data[18,2] = [double]0.1; $data[18,3] = [double]0.4; $data[18,4] = [double]513.6
$data[19,0] = "2013-07"; $data[19,1] = [double]-100; $data[19,2] = [double]0; $data[19,3] = [double]0.3; $data[19,4] = [double]236.4
$data[20,0] = "2013-09"; $data[20,1] = [double]647.3; $data[20,2] = [double]0.2; $data[20,3] = [double]0.9; $data[20,4] = [double]140.1
$data[21,0] = "2014-10"; $data[21,1] = ""; $data[21,2] = ""; $data[21,3] = [double]0; $data[21,4] = ""
$data[22,0] = "2014-11"; $data[22,1] = ""; $data[22,2] = ""; $data[22,3] = [double]0; $data[22,4] = ""
$data[23,0] = "2014-12"; $data[23,1] = ""; $data[23,2] = ""; $data[23,3] = [double]0; $data[23,4] = ""
$data[24,0] = "2014-02"; $data[24,1] = ""; $data[24,2] = ""; $data[24,3] = ""; $data[24,4] = [double]-49.9
$data[25,0] = "2014-04"; $data[25,1] = ""; $data[25,2] = ""; $data[25,3] = [double]0; $data[25,4] = ""
$data[26,0] = "2014-05"; $data[26,1] = ""; $data[26,2] = ""; $data[26,3] = [double]0; $data[26,4] = ""
$data[27,0] = "2014-07"; $data[27,1] = ""; $data[27,2] = ""; $data[27,3] = [double]1.4; $data[27,4] = ""
$data[28,0] = "2014-08"; $data[28,1] = ""; $data[28,2] = ""; $data[28,3] = [double]0; $data[28,4] = ""
$data[29,0] = "2014-09"; $data[29,1] = ""; $data[29,2] = ""; $data[29,3] = [double]0; $data[29,4] = ""
$data[30,0] = "2015-10"; $data[30,1] = [double]-100; $data[30,2] = [double]0.2; $data[30,3] = [double]3.6; $data[30,4] = [double]-100
$data[31,0] = "2015-11"; $data[31,1] = ""; $data[31,2] = ""; $data[31,3] = [double]3.8; $data[31,4] = ""
$data[32,0] = "2015-12"; $data[32,1] = ""; $data[32,2] = ""; $data[32,3] = [double]5.4; $data[32,4] = ""
$data[33,0] = "2015-02"; $data[33,1] = ""; $data[33,2] = ""; $data[33,3] = ""; $data[33,4] = [double]-100
$data[34,0] = "2015-03"; $data[34,1] = ""; $data[34,2] = ""; $data[34,3] = [double]1.3; $data[34,4] = ""
$data[35,0] = "2015-04"; $data[35,1] = ""; $data[35,2] = ""; $data[35,3] = [double]2.1; $data[35,4] = ""
$data[36,0] = "2015-05"; $data[36,1] = ""; $data[36,2] = ""; $data[36,3] = [double]2.4; $data[36,4] = ""
$data[37,0] = "2015-06"; $data[37,1] = ""; $data[37,2] = ""; $data[37,3] = [double]2.7; $data[37,4] = ""
$data[38,0] = "2015-07"; $data[38,1] = ""; $data[38,2] = ""; $data[38,3] = [double]3; $data[38,4] = ""
$data[39,0] = "2015-08"; $data[39,1] = ""; $data[39,2] = ""; $data[39,3] = [double]3.2; $data[39,4] = ""
$data[40,0] = "2015-09"; $data[40,1] = ""; $data[40,2] = ""; $data[40,3] = [double]3.4; $data[40,4] = ""
$data[41,0] = "2016-10"; $data[41,1] = [double]-100; $data[41,2] = [double]0; $data[41,3] = [double]0.7; $data[41,4] = [double]-75
$data[42,0] = "2016-11"; $data[42,1] = [double]-100; $data[42,2] = [double]0; $data[42,3] = [double]0.3; $data[42,4] = [double]-90
$data[43,0] = "2016-12"; $data[43,1] = [double]-100; $data[43,2] = [double]0; $data[43,3] = [double]0.3; $data[43,4] = [double]-93.5
$data[44,0] = "2016-03"; $data[44,1] = ""; $data[44,2] = ""; $data[44,3] = ""; $data[44,4] = [double]-90
$data[45,0] = "2016-04"; $data[45,1] = ""; $data[45,2] = ""; $data[45,3] = ""; $data[45,4] = [double]-94.09999999999999
$data[46,0] = "2016-05"; $data[46,1] = [double]-100; $data[46,2] = [double]0; $data[46,3] = [double]0.1; $data[46,4] = [double]-94.40000000000001
$data[47,0] = "2016-06"; $data[47,1] = [double]-100; $data[47,2] = [double]0; $data[47,3] = [double]0.1; $data[47,4] = [double]-95
$data[48,0] = "2016-07"; $data[48,1] = [double]-100; $data[48,2] = [double]0; $data[48,3] = [double]0.2; $data[48,4] = [double]-90.90000000000001
$data[49,0] = "2016-08"; $data[49,1] = [double]-100; $data[49,2] = [double]0; $data[49,3] = [double]0.2; $data[49,4] = [double]-91.7
$data[50,0] = "2016-09"; $data[50,1] = [double]150; $data[50,2] = [double]0.5; $data[50,3] = [double]0.7; $data[50,4] = [double]-73.09999999999999
$data[51,0] = "2017-10"; $data[51,1] = [double]0; $data[51,2] = [double]0; $data[51,3] = [double]0.4; $data[51,4] = [double]-42.9
$data[52,0] = "2017-11"; $data[52,1] = [double]0; $data[52,2] = [double]0; $data[52,3] = [double]0.4; $data[52,4] = [double]100
$data[53,0] = "2017-12"; $data[53,1] = [double]0; $data[53,2] = [double]0; $data[53,3] = [double]0.5; $data[53,4] = [double]150
$data[54,0] = "2017-02"; $data[54,1] = ""; $data[54,2] = ""; $data[54,3] = ""; $data[54,4] = [double]0
$data[55,0] = "2017-03"; $data[55,1] = [double]0; $data[55,2] = [double]0.2; $data[55,3] = [double]0.5; $data[55,4] = [double]400
$data[56,0] = "2017-04"; $data[56,1] = [double]0; $data[56,2] = [double]0.2; $data[56,3] = [double]0.8; $data[56,4] = [double]700
$data[57,0] = "2017-05"; $data[57,1] = [double]0; $data[57,2] = [double]0.1; $data[57,3] = [double]0.7; $data[57,4] = [double]600
$data[58,0] = "2017-06"; $data[58,1] = [double]0; $data[58,2] = [double]0.1; $data[58,3] = [double]0.8; $data[58,4] = [double]700
$data[59,0] = "2017-07"; $data[59,1] = [double]0; $data[59,2] = [double]0; $data[59,3] = [double]0.9; $data[59,4] = [double]350
$data[60,0] = "2017-08"; $data[60,1] = [double]0; $data[60,2] = [double]0; $data[60,3] = [double]0.4; $data[60,4] = [double]100
$data[61,0] = "2017-09"; $data[61,1] = [double]-100; $data[61,2] = [double]0; $data[61,3] = [double]0.4; $data[61,4] = [double]-42.9

# Write the reordered values back (header row A1:E1 is untouched).
$ws.Range("A2:E63").Value = $data

# Drop column F entirely (its header and all 62 data values).
$ws.Columns("F").Delete()

